$d = $word.ActiveDocument

# --- Step 1: move the "_GoBack" bookmark from the end of paragraph 3 to
# the end of paragraph 1 (right after "16/7/18", still inside that paragraph).
# We insert a temporary character, wrap the bookmark around it (a
# non-collapsed range), then delete the temporary character again. Doing
# it this way leaves both bookmarkStart/bookmarkEnd sitting cleanly right
# after the existing run instead of straddling the paragraph boundary.
$p1 = $d.Paragraphs(1).Range
$p1.MoveEnd(1, -1)
$p1.Collapse(0)
$p1.InsertAfter("X")
$tmp = $d.Range($p1.Start, $p1.Start + 1)
$d.Bookmarks.Add("_GoBack", $tmp)
$tmp.Delete()

# --- Step 2: remove paragraph 3 ("Update GM to set attack/defender to
# arrays") entirely -- this also removes the old bookmark that used to
# live there.
$d.Paragraphs(3).Range.Delete()

# --- Step 3: rewrite paragraph 2's text.
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("Update energy Lines to get deleted when you right click on the, might need to add collision box", $false, $false, $false, $false, $false, $true, 1, $false, "Bugs- 2+ lines (needs testing)", 2)

# Add the second sentence as its own run (not merged into the first) by
# inserting it as a new paragraph and then merging the paragraph mark
# back out -- the two pieces of text stay as separate <w:r> runs.
$p2end = $d.Paragraphs(2).Range
$p2end.Collapse(0)
$p2end.InsertParagraphAfter()

$p3 = $d.Paragraphs(3).Range
$p3.InsertAfter(", sharing energy, setting target to self")

$p2merge = $d.Paragraphs(2).Range
$p2merge.MoveEnd(1, -1)
$markPos = $p2merge.End
$markRange = $d.Range($markPos, $markPos + 1)
$markRange.Delete()
